# Updating hyperparameters table: four new (better) training runs were
# logged at the top of the results table, pushing the previous rows down.
# The sheet is kept sorted ascending by Validation Loss (column F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the four new result rows above the existing data (old
# row 2 becomes row 6, ..., old row 8 becomes row 12).
$ws.Rows("2:5").Insert()

# The inserted rows pick up row 1's formatting (header style) by default;
# the new data rows should be unformatted like the rest of the table.
$ws.Range("A2:F5").ClearFormats()

# New hyperparameter runs: Iterations, Learning Rate, Hidden Nodes,
# Output Nodes, Training loss, Validation Loss.
$newRuns = @(
    @(11000, 0.3, 8,  1, 0.058, 0.134),
    @(11000, 0.3, 7,  1, 0.059, 0.141),
    @(11000, 0.3, 10, 1, 0.058, 0.144),
    @(10000, 0.3, 6,  1, 0.062, 0.15)
)

$r = 2
foreach ($run in $newRuns) {
    $ws.Cells.Item($r, 1).Value = $run[0]
    $ws.Cells.Item($r, 2).Value = $run[1]
    $ws.Cells.Item($r, 3).Value = $run[2]
    $ws.Cells.Item($r, 4).Value = $run[3]
    $ws.Cells.Item($r, 5).Value = $run[4]
    $ws.Cells.Item($r, 6).Value = $run[5]
    $r++
}

# Keep the table sorted ascending by Validation Loss, as before.
$sortRange = $ws.Range("A2:F12")
$sortKey = $ws.Range("F1")
$sortRange.Sort($sortKey, 1) | Out-Null

# Match the author's final selection.
$ws.Range("E3").Select() | Out-Null
